$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-07-09 Tuesday" "2024-07-10 Wednesday"

Replace-Text "90×50=" "73×54="
Replace-Text "60×69=" "23×29="
Replace-Text "95×65=" "98×85="
Replace-Text "25×87=" "56×98="
Replace-Text "62×51=" "26×26="
Replace-Text "96×61=" "87×85="
Replace-Text "25×15=" "55×35="
Replace-Text "90×98=" "76×33="
Replace-Text "73×24=" "65×99="
Replace-Text "75×24=" "51×73="
Replace-Text "50×73=" "16×28="
Replace-Text "27×65=" "25×77="
Replace-Text "65×28=" "57×40="
Replace-Text "11×43=" "84×34="
Replace-Text "19×68=" "78×67="
Replace-Text "67×49=" "25×70="
Replace-Text "73×43=" "60×57="
Replace-Text "86×78=" "36×26="
Replace-Text "21×71=" "29×11="
Replace-Text "31×56=" "57×24="
Replace-Text "33×73=" "61×21="
Replace-Text "45×78=" "42×86="
Replace-Text "97×91=" "14×63="
Replace-Text "47×26=" "95×79="
Replace-Text "60×36=" "95×83="
